# Update Name of Algo
# Apply updated imputed values to the result_data_KNN worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.068
$ws.Range("D3").Value = -7.507
$ws.Range("A12").Value = -21.519
$ws.Range("B14").Value = 5.944
$ws.Range("B26").Value = 6.225
$ws.Range("D30").Value = -7.257
$ws.Range("B31").Value = 6.805
$ws.Range("A32").Value = -21.351
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -20.945
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.093
$ws.Range("D44").Value = -7.874
$ws.Range("B45").Value = 5.75
$ws.Range("A46").Value = -21.553
$ws.Range("A54").Value = -22.209
$ws.Range("A55").Value = -22.21
$ws.Range("B57").Value = 5.331999999999999
$ws.Range("D58").Value = -8.303000000000001
$ws.Range("A67").Value = -21.586
$ws.Range("A69").Value = -21.637
$ws.Range("A72").Value = -21.55
$ws.Range("D84").Value = -8.260999999999999
$ws.Range("D89").Value = -6.962000000000001
$ws.Range("A91").Value = -21.522
$ws.Range("D91").Value = -6.931999999999999
$ws.Range("A99").Value = -20.828
$ws.Range("B100").Value = 5.558
$ws.Range("B102").Value = 7.468000000000001
$ws.Range("D102").Value = -7.873
